$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 74
$ws1.Range("F4").Value = 2078
$ws1.Range("F6").Value = 624
$ws1.Range("F7").Value = 101
$ws1.Range("F9").Value = 10676
$ws1.Range("F15").Value = 7538
$ws1.Range("F16").Value = 1113
$ws1.Range("F17").Value = 719
$ws1.Range("F18").Value = 257

# Sheet "全部类型" (sheet4): update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 74
$ws4.Range("F4").Value = 2078
$ws4.Range("F6").Value = 624
$ws4.Range("F8").Value = 101
$ws4.Range("F12").Value = 10676
$ws4.Range("F18").Value = 7538
$ws4.Range("F19").Value = 1113
$ws4.Range("F20").Value = 719
$ws4.Range("F21").Value = 257
